# Horarios actualizados Linea 141 - 583
# Applies the 12:11:21 scrape refresh across all three sheets:
#   - LP1912      : rows 62/63 re-ordered + rows 126-158 rewritten (6 new rows appended, data re-sorted by arrival time)
#   - LP1912-215  : only the "Ultima actualizacion" timestamp changes
#   - 6203-6173   : only the "Ultima actualizacion" timestamp changes + one new row (29) appended

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 12:11:21"
$ws1.Range("A3").Value = "Total filas: 153"

$rows1 = @(
  @(62, "07:49:32", "09:17", "14_ABASTO", 88, "LP1912"),
  @(63, "08:38:24", "09:17", "27_EL RETIRO", 39, "LP1912"),
  @(126, "12:11:21", "12:11", "16_P MOR-SANTA ANA", 0, "LP1912"),
  @(127, "12:11:21", "12:12", "15_ABASTO", 1, "LP1912"),
  @(128, "11:33:52", "12:16", "10_OLMOS", 43, "LP1912"),
  @(129, "11:13:15", "12:17", "10_OLMOS", 64, "LP1912"),
  @(130, "10:36:50", "12:21", "215C_EL PATO", 105, "LP1912"),
  @(131, "10:56:15", "12:22", "215C_EL PATO", 86, "LP1912"),
  @(132, "11:13:15", "12:29", "23_HERNANDEZ", 76, "LP1912"),
  @(133, "10:36:50", "12:32", "14_ABASTO", 116, "LP1912"),
  @(134, "11:33:52", "12:32", "23_HERNANDEZ", 59, "LP1912"),
  @(135, "10:56:15", "12:33", "14_ABASTO", 97, "LP1912"),
  @(136, "10:56:15", "12:33", "27_EL RETIRO", 97, "LP1912"),
  @(137, "10:36:50", "12:34", "15_ABASTO", 118, "LP1912"),
  @(138, "11:46:32", "12:34", "23_HERNANDEZ", 48, "LP1912"),
  @(139, "10:49:38", "12:36", "27_EL RETIRO", 107, "LP1912"),
  @(140, "11:53:44", "12:36", "23_HERNANDEZ", 43, "LP1912"),
  @(141, "12:11:21", "12:37", "23_HERNANDEZ", 26, "LP1912"),
  @(142, "12:11:21", "12:37", "27_EL RETIRO", 26, "LP1912"),
  @(143, "11:33:52", "12:47", "14_ABASTO", 74, "LP1912"),
  @(144, "11:33:52", "12:48", "15X38_ABASTO", 75, "LP1912"),
  @(145, "10:49:38", "12:48", "16_SANTA ANA", 119, "LP1912"),
  @(146, "11:33:52", "13:02", "11_ETCHEVERRY", 89, "LP1912"),
  @(147, "11:13:15", "13:03", "11_ETCHEVERRY", 110, "LP1912"),
  @(148, "11:33:52", "13:03", "215C_EL PATO", 90, "LP1912"),
  @(149, "11:46:32", "13:04", "215C_EL PATO", 78, "LP1912"),
  @(150, "11:33:52", "13:13", "16_SANTA ANA", 100, "LP1912"),
  @(151, "11:33:52", "13:17", "10_OLMOS", 104, "LP1912"),
  @(152, "11:53:44", "13:21", "23_HERNANDEZ", 88, "LP1912"),
  @(153, "12:11:21", "13:24", "23_HERNANDEZ", 73, "LP1912"),
  @(154, "11:33:52", "13:25", "16_P MOR-SANTA ANA", 112, "LP1912"),
  @(155, "11:53:44", "13:32", "215A_EL PATO", 99, "LP1912"),
  @(156, "12:11:21", "13:32", "14_ABASTO", 81, "LP1912"),
  @(157, "11:46:32", "13:33", "215A_EL PATO", 107, "LP1912"),
  @(158, "11:53:44", "13:47", "225_GOMEZ", 114, "LP1912")
)

foreach ($r in $rows1) {
    $ws1.Cells.Item($r[0], 1).Value = $r[1]
    $ws1.Cells.Item($r[0], 2).Value = $r[2]
    $ws1.Cells.Item($r[0], 3).Value = $r[3]
    $ws1.Cells.Item($r[0], 4).Value = $r[4]
    $ws1.Cells.Item($r[0], 5).Value = $r[5]
}

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215 (only the refresh timestamp changes)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 12:11:21"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173 (timestamp + one new appended row)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 12:11:21"
$ws3.Range("A3").Value = "Total filas: 24"

$ws3.Cells.Item(29, 1).Value = "12:11:21"
$ws3.Cells.Item(29, 2).Value = "13:57"
$ws3.Cells.Item(29, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(29, 4).Value = 106
$ws3.Cells.Item(29, 5).Value = "L6203"
